$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows 14-16 down to 15-17.
# Excel.Insert() copies formatting (e.g. date style on column D) from the row above.
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new record's data.
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C14").Value = "Ñuble"
$ws.Range("D14").Value = 44524
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = 100112040
$ws.Range("G14").Value = "Cilantro"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 2000
$ws.Range("N14").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O14").Value = "Provincia de Diguillín"
$ws.Range("P14").Value = 2000
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
